$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 316, shifting existing rows 316:346 down to 317:347
$ws.Rows.Item(316).Insert()

# Populate the newly inserted row 316 with data (copy of surrounding record, with
# the fields below updated to the new observation's values)
$ws.Cells.Item(316, 1).Value = 9
$ws.Cells.Item(316, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(316, 3).Value = "Metropolitana"
$ws.Cells.Item(316, 4).Value = 44858
$ws.Cells.Item(316, 5).Value = 13
$ws.Cells.Item(316, 6).Value = 100112021
$ws.Cells.Item(316, 7).Value = "Ají"
$ws.Cells.Item(316, 8).Value = "Americana (o)"
$ws.Cells.Item(316, 9).Value = "Primera"
$ws.Cells.Item(316, 10).Value = 120
$ws.Cells.Item(316, 11).Value = 2000
$ws.Cells.Item(316, 12).Value = 2000
$ws.Cells.Item(316, 13).Value = 2000
$ws.Cells.Item(316, 14).Value = "$/kilo"
$ws.Cells.Item(316, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(316, 16).Value = 2000
$ws.Cells.Item(316, 17).Value = 1
$ws.Cells.Item(316, 18).Value = "Hortaliza"
